$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Source Control" / "Set up git" story entirely (row 2)
$ws.Rows(2).Delete()

# Insert a blank row for the new "Links to corresponding videos..." story (now row 3)
$ws.Rows(3).Insert()

# Insert a blank row for the new "Investigate embedding YouTube videos..." story (now row 6)
$ws.Rows(6).Insert()

# Fill in the text for the stories that moved into newly-inserted rows / got reworded
$ws.Range("B3").Value = "Links to corresponding videos on YouTube or our site (with embedded YouTube videos) beside songs."
$ws.Range("B5").Value = "Put individual links to videos on YouTube on Video pages"
$ws.Range("B6").Value = "Investigate embedding YouTube videos on pages"
$ws.Range("B7").Value = "Remove extra space on right hand side of pages (centre content)"

# Add the new stories at the bottom of the backlog
$ws.Range("A10").Value = "Site"
$ws.Range("B10").Value = "Create SPA using AngularJS"
$ws.Range("A11").Value = "Social Links"
$ws.Range("B11").Value = "Twitter"
$ws.Range("B12").Value = "Facebook Like"
$ws.Range("B13").Value = "Facebook Comment"
$ws.Range("B14").Value = "Facebook Share"
$ws.Range("B15").Value = "Others"

# Reword the zip files story
$ws.Range("B4").Value = "Use real zip files rather than ones created on the fly or make sure these don't fail to download"

# Update selection to match the author's saved cursor position
$ws.Range("B5").Select()
